# #5: property aircraft done
# The "航空器" (Aircraft) property sheet was an unfilled template (no real
# data rows) and is removed entirely. The "category" columns on the
# "建物" (building) and "汽車" (car) sheets - which had been left as the
# generic "land" placeholder value copied from the land-sheet template -
# are corrected to "building" and "car" respectively.

$wb = $excel.ActiveWorkbook

# Delete the now-unused "航空器" (Aircraft) placeholder sheet.
$aircraft = $wb.Worksheets.Item("航空器")
$aircraft.Delete()

# Fix the "building" sheet's category column (I2:I11) from "land" to "building".
$buildingSheet = $wb.Worksheets.Item("建物")
$buildingRange = $buildingSheet.Range("I2:I11")
for ($r = 1; $r -le $buildingRange.Rows.Count; $r++) {
    $buildingRange.Cells.Item($r, 1).Value = "building"
}

# Fix the "car" sheet's category column (H2:H3) from "land" to "car".
$carSheet = $wb.Worksheets.Item("汽車")
$carRange = $carSheet.Range("H2:H3")
for ($r = 1; $r -le $carRange.Rows.Count; $r++) {
    $carRange.Cells.Item($r, 1).Value = "car"
}
